# ------------------------------------------------------------------
# Workers_problem.xlsx edit script
#  1. Rename existing sheets: Sheet1 -> Workers1, Sheet2 -> Worker2
#  2. Tweak a few values on Workers1 (D4/E4 + derived SUMPRODUCT totals)
#     and add the new "Decision Variable" vector label cell (D2)
#  3. Add a brand-new sheet named "Sheet1" (TV-advertising LP model),
#     make it the active sheet/tab
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. rename existing sheets -----------------------------------
$wsWorkers1 = $wb.Worksheets.Item(1)
$wsWorker2  = $wb.Worksheets.Item(2)
$wsWorkers1.Name = "Workers1"
$wsWorker2.Name  = "Worker2"

# ---- 2. small data tweaks on Workers1 ----------------------------
$wsWorkers1.Range("D2").Value = "[18.0, 10.0, 8.0, 0.0, 13.0, 0.0, 4.0]"

$wsWorkers1.Range("D4").Value = 18
$wsWorkers1.Range("E4").Value = 10
# F4,G4,H4,I4,J4 stay 8,0,13,0,4 (unchanged)

# L8 / L13 totals recompute automatically from the SUMPRODUCT formula
# already in place once D4/E4 change, but set explicitly to be safe.
$wsWorkers1.Calculate()

# ---- 3. add the new "Sheet1" (Advertising Mode LP) ---------------
$wsAdv = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsWorker2)
$wsAdv.Name = "Sheet1"

# Title
$wsAdv.Range("A1").Value = "Advertising Mode"

# ---- header row (5): Group / show names / "Minimum Exposures" ----
$wsAdv.Range("B5").Value = "Group"
$wsAdv.Range("C5").Value = "Revenge"
$wsAdv.Range("D5").Value = "Sunday Night Football"
$wsAdv.Range("E5").Value = "The Simpsons"
$wsAdv.Range("F5").Value = "Sports Center"
$wsAdv.Range("G5").Value = "Homeland"
$wsAdv.Range("H5").Value = "Rachel Ray"
$wsAdv.Range("I5").Value = "CNN"
$wsAdv.Range("J5").Value = "The Good Wife"
$wsAdv.Range("L5").Value = "Minimum Exposures"

# ---- demographic rows 6-11 ----------------------------------------
$wsAdv.Range("B6").Value  = "Male 18-35"
$wsAdv.Range("C6:J6").Value = @(5, 6, 5, 0.5, 0.7, 0.1, 0.1, 3)
$wsAdv.Range("L6").Value  = 60

$wsAdv.Range("B7").Value  = "Male 36-55"
$wsAdv.Range("C7:J7").Value = @(3, 5, 2, 0.5, 0.2, 0.1, 0.2, 5)
$wsAdv.Range("L7").Value  = 60

$wsAdv.Range("B8").Value  = "Male > 55"
$wsAdv.Range("C8:J8").Value = @(1, 3, 0, 0.3, 0, 0, 0.3, 4)
$wsAdv.Range("L8").Value  = 28

$wsAdv.Range("B9").Value  = "Female 18-35"
$wsAdv.Range("C9:J9").Value = @(6, 1, 4, 0.1, 0.9, 0.6, 0.1, 3)
$wsAdv.Range("L9").Value  = 60

$wsAdv.Range("B10").Value = "Female 36-55"
$wsAdv.Range("C10:J10").Value = @(4, 1, 2, 0.1, 0.1, 1.3, 0.2, 5)
$wsAdv.Range("L10").Value = 60

$wsAdv.Range("B11").Value = "Female > 55"
$wsAdv.Range("C11:J11").Value = @(2, 1, 0, 0, 0, 0.4, 0.3, 4)
$wsAdv.Range("L11").Value = 28

# ---- cost per ad row 12 --------------------------------------------
$wsAdv.Range("B12").Value = "Cost Per Ad"
$wsAdv.Range("C12:J12").Value = @(140, 100, 80, 9, 13, 15, 8, 140)

# ---- second "styled" header block (row 15) + input row (16) -------
$wsAdv.Range("B15").Value = "Decision Variable"
$wsAdv.Range("C15").Value = "Revenge"
$wsAdv.Range("D15").Value = "Sunday Night Football"
$wsAdv.Range("E15").Value = "The Simpsons"
$wsAdv.Range("F15").Value = "Sports Center"
$wsAdv.Range("G15").Value = "Homeland"
$wsAdv.Range("H15").Value = "Rachel Ray"
$wsAdv.Range("I15").Value = "CNN"
$wsAdv.Range("J15").Value = "The Good Wife"
$wsAdv.Range("L15").Value = "Obj"

$wsAdv.Range("B16").Value = "Number ads purchased"
$wsAdv.Range("L16").Formula = "=SUMPRODUCT(C12:J12,C16:J16)"

# ---- constraints block (rows 19-26) --------------------------------
$wsAdv.Range("B19").Value = "Constraints"
$wsAdv.Range("B20").Value = "Group"

$wsAdv.Range("B21").Value = "Male 18-35"
$wsAdv.Range("C21").Formula = "=SUMPRODUCT($C$16:$J$16,C6:J6)"
$wsAdv.Range("D21").Value = ">="
$wsAdv.Range("E21").FormulaArray = "=L6:L11"

$wsAdv.Range("B22").Value = "Male 36-55"
$wsAdv.Range("C22").Formula = "=SUMPRODUCT($C$16:$J$16,C7:J7)"
$wsAdv.Range("D22").Value = ">="
$wsAdv.Range("E22").Value = 60

$wsAdv.Range("B23").Value = "Male > 55"
$wsAdv.Range("C23").Formula = "=SUMPRODUCT($C$16:$J$16,C8:J8)"
$wsAdv.Range("D23").Value = ">="
$wsAdv.Range("E23").Value = 28

$wsAdv.Range("B24").Value = "Female 18-35"
$wsAdv.Range("C24").Formula = "=SUMPRODUCT($C$16:$J$16,C9:J9)"
$wsAdv.Range("D24").Value = ">="
$wsAdv.Range("E24").Value = 60

$wsAdv.Range("B25").Value = "Female 36-55"
$wsAdv.Range("C25").Formula = "=SUMPRODUCT($C$16:$J$16,C10:J10)"
$wsAdv.Range("D25").Value = ">="
$wsAdv.Range("E25").Value = 60

$wsAdv.Range("B26").Value = "Female > 55"
$wsAdv.Range("C26").Formula = "=SUMPRODUCT($C$16:$J$16,C11:J11)"
$wsAdv.Range("D26").Value = ">="
$wsAdv.Range("E26").Value = 28

# ---- column widths (visual match) ----------------------------------
$wsAdv.Columns.Item("B").ColumnWidth = 15.7109375
$wsAdv.Columns.Item("C").ColumnWidth = 11
$wsAdv.Columns.Item("D").ColumnWidth = 23.28515625
$wsAdv.Columns.Item("E").ColumnWidth = 15.85546875
$wsAdv.Columns.Item("F").ColumnWidth = 15.7109375
$wsAdv.Columns.Item("G").ColumnWidth = 12.5703125
$wsAdv.Columns.Item("H").ColumnWidth = 13.140625
$wsAdv.Columns.Item("J").ColumnWidth = 16.85546875
$wsAdv.Columns.Item("K").ColumnWidth = 8.85546875
$wsAdv.Columns.Item("L").ColumnWidth = 14.85546875

# ---- styling: blue/bold header (row15 C:J) + left-border header (B20) ----
# OLE (BGR-packed) colour integers corresponding to the theme colours used
# by the authored workbook: accent1 solid (4472C4), white (FFFFFF),
# accent1 25% (light fill, DAE3F3) and accent1 40%-tint (border, 8FAADC).
$colBlueSolid  = 12874308   # RGB 4472C4 - "blue, accent1"
$colWhite      = 16777215   # RGB FFFFFF - white font on the blue header
$colLightFill  = 15983578   # RGB DAE3F3 - "blue, accent1, lighter 80%"
$colBorder     = 14461583   # RGB 8FAADC - "blue, accent1, lighter 40%"

function Style-HeaderCell($rng, [bool]$rightEdge) {
    $rng.Font.Bold = $true
    $rng.Font.Color = $colWhite
    $rng.Interior.Color = $colBlueSolid
    $rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders.Item(8).Weight = 2      # xlThin
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = $colBorder
    if ($rightEdge) {
        $rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
        $rng.Borders.Item(10).Weight = 2
        $rng.Borders.Item(10).Color = $colBorder
    }
}

function Style-InputCell($rng, [bool]$rightEdge) {
    $rng.Interior.Color = $colLightFill
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = $colBorder
    if ($rightEdge) {
        $rng.Borders.Item(10).LineStyle = 1
        $rng.Borders.Item(10).Weight = 2
        $rng.Borders.Item(10).Color = $colBorder
    }
}

function Style-LeftHeaderCell($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Color = $colWhite
    $rng.Interior.Color = $colBlueSolid
    $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(7).Color = $colBorder
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = $colBorder
}

function Style-LeftFillCell($rng) {
    $rng.Interior.Color = $colLightFill
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(7).Color = $colBorder
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = $colBorder
}

function Style-LeftPlainCell($rng) {
    $rng.Borders.Item(7).LineStyle = 1
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(7).Color = $colBorder
    $rng.Borders.Item(8).LineStyle = 1
    $rng.Borders.Item(8).Weight = 2
    $rng.Borders.Item(8).Color = $colBorder
    $rng.Borders.Item(9).LineStyle = 1
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(9).Color = $colBorder
}

# row 15 header (blue/bold), C:I then J gets the extra right edge
Style-HeaderCell $wsAdv.Range("C15:I15") $false
Style-HeaderCell $wsAdv.Range("J15") $true

# row 16 input row (light fill), C:I then J gets the extra right edge
Style-InputCell $wsAdv.Range("C16:I16") $false
Style-InputCell $wsAdv.Range("J16") $true

# left-edge styled column B for the small constraints report (rows 20-26)
Style-LeftHeaderCell $wsAdv.Range("B20")
Style-LeftFillCell   $wsAdv.Range("B21")
Style-LeftPlainCell  $wsAdv.Range("B22")
Style-LeftFillCell   $wsAdv.Range("B23")
Style-LeftPlainCell  $wsAdv.Range("B24")
Style-LeftFillCell   $wsAdv.Range("B25")
Style-LeftPlainCell  $wsAdv.Range("B26")

# ---- build the Table3 ListObject over B5:J12 -----------------------
$loAdv = $wsAdv.ListObjects.Add(1, $wsAdv.Range("B5:J12"), [System.Reflection.Missing]::Value, 1)
$loAdv.Name = "Table3"

# ---- view settings ---------------------------------------------------
$wsAdv.Application.ActiveWindow.ScrollRow = 1
$wsAdv.Range("K4").Select()

$wsWorkers1.Application.ActiveWindow.ScrollRow = 1

$wb.Windows.Item(1).Activate()
$wsAdv.Activate()

$wb.Save()
